$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 1.08
$ws.Range("N4").Value = 1.26
$ws.Range("O4").Value = 3.15
$ws.Range("R4").Value = 2.05
$ws.Range("S4").Value = 1.6
$ws.Range("T4").Value = 120
$ws.Range("U4").Value = 500
$ws.Range("V4").Value = 110
$ws.Range("X4").Value = 700
$ws.Range("Y4").Value = 300
$ws.Range("Z4").Value = 25
$ws.Range("AA4").Value = 22
$ws.Range("AC4").Value = 150
$ws.Range("AD4").Value = 13
$ws.Range("AE4").Value = 7.8
$ws.Range("AF4").Value = 13
$ws.Range("AG4").Value = 6.8
$ws.Range("AI4").Value = 35
$ws.Range("AJ4").Value = 1000

# Row 5
$ws.Range("G5").Value = 6.6
$ws.Range("H5").Value = 4.65
$ws.Range("I5").Value = 1.38
$ws.Range("N5").Value = 1.5
$ws.Range("O5").Value = 2.25
$ws.Range("R5").Value = 1.72
$ws.Range("S5").Value = 1.9
$ws.Range("T5").Value = 22
$ws.Range("U5").Value = 45
$ws.Range("X5").Value = 65
$ws.Range("Z5").Value = 15.5
$ws.Range("AA5").Value = 9.5
$ws.Range("AC5").Value = 65
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 7.5
$ws.Range("AH5").Value = 10.5
$ws.Range("AI5").Value = 22
$ws.Range("AJ5").Value = 450

# Row 6
$ws.Range("G6").Value = 1.17
$ws.Range("H6").Value = 6.5
$ws.Range("I6").Value = 13
$ws.Range("N6").Value = 1.31
$ws.Range("O6").Value = 2.87
$ws.Range("R6").Value = 1.85
$ws.Range("S6").Value = 1.75
$ws.Range("T6").Value = 11
$ws.Range("U6").Value = 7.5
$ws.Range("V6").Value = 10.25
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 10
$ws.Range("Y6").Value = 26
$ws.Range("Z6").Value = 22
$ws.Range("AA6").Value = 15
$ws.Range("AC6").Value = 90
$ws.Range("AD6").Value = 50
$ws.Range("AE6").Value = 120
$ws.Range("AG6").Value = 500
$ws.Range("AI6").Value = 100
$ws.Range("AJ6").Value = 600

# Row 11
$ws.Range("G11").Value = 1.91
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 3.1
$ws.Range("N11").Value = 1.57
$ws.Range("O11").Value = 2.35
$ws.Range("T11").Value = 11
$ws.Range("U11").Value = 12
$ws.Range("W11").Value = 19
$ws.Range("Y11").Value = 21
$ws.Range("Z11").Value = 17
$ws.Range("AA11").Value = 8
$ws.Range("AB11").Value = 13
$ws.Range("AC11").Value = 41
$ws.Range("AD11").Value = 15
$ws.Range("AF11").Value = 12
$ws.Range("AG11").Value = 34
$ws.Range("AH11").Value = 23
$ws.Range("AI11").Value = 26
$ws.Range("AJ11").Value = 126

# Row 16
$ws.Range("N16").Value = 1.77
$ws.Range("O16").Value = 1.92

# Row 18
$ws.Range("J18").Value = 1.06
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 1.36
$ws.Range("M18").Value = 3
$ws.Range("N18").Value = 2.1
$ws.Range("O18").Value = 1.7
$ws.Range("AJ18").Value = 501

# Row 30
$ws.Range("K30").Value = 17

# Row 31
$ws.Range("J31").Value = 1.05
$ws.Range("K31").Value = 11
$ws.Range("N31").Value = 1.83
$ws.Range("O31").Value = 1.98

# Row 33
$ws.Range("L33").Value = 1.18
$ws.Range("M33").Value = 4.5
$ws.Range("N33").Value = 1.62
$ws.Range("O33").Value = 2.25

# Row 34
$ws.Range("G34").Value = 1.48
$ws.Range("I34").Value = 5.5
$ws.Range("N34").Value = 1.36
$ws.Range("O34").Value = 3.1
$ws.Range("R34").Value = 1.44
$ws.Range("S34").Value = 2.63
$ws.Range("AB34").Value = 13
$ws.Range("AD34").Value = 23
$ws.Range("AE34").Value = 34
$ws.Range("AF34").Value = 17
$ws.Range("AH34").Value = 34
$ws.Range("AI34").Value = 29

# Row 37
$ws.Range("H37").Value = 3.1
$ws.Range("I37").Value = 3.05
$ws.Range("L37").Value = 1.29
$ws.Range("M37").Value = 2.95
$ws.Range("N37").Value = 1.87
$ws.Range("O37").Value = 1.75
$ws.Range("P37").Value = 1.42
$ws.Range("Q37").Value = 2.47
$ws.Range("R37").Value = 1.65
$ws.Range("S37").Value = 1.98
$ws.Range("T37").Value = 8.25
$ws.Range("U37").Value = 11.75
$ws.Range("V37").Value = 8.75
$ws.Range("W37").Value = 24
$ws.Range("X37").Value = 18
$ws.Range("Y37").Value = 26
$ws.Range("Z37").Value = 9.25
$ws.Range("AA37").Value = 6.1
$ws.Range("AB37").Value = 13
$ws.Range("AC37").Value = 60
$ws.Range("AD37").Value = 9.25
$ws.Range("AE37").Value = 16
$ws.Range("AF37").Value = 10.75
$ws.Range("AG37").Value = 40
$ws.Range("AH37").Value = 28
$ws.Range("AJ37").Value = 450
